# Actor.xlsx -> ActorTable: add attackSpeed|Float column (E) ------------------
$wb = $excel.ActiveWorkbook

$actorTable = $wb.Worksheets.Item("ActorTable")

# New header in E1 (this also creates the "attackSpeed|Float" shared string)
$actorTable.Cells.Item(1, 5).Value = "attackSpeed|Float"

# Fill E2:E18 with 1, matching the rest of the table's attackDelay-style flag column
$actorTable.Range("E2:E18").Value = 1

# Match the new column's width to the rest of the stat columns (target ~16.25)
$actorTable.Columns.Item(5).ColumnWidth = 15.5

# ------------------------------------------------------------------------
# View-state bookkeeping to mirror the saved workbook:
#  - ActorPowerLevelTable's selection moves to F1 (first empty column there)
#  - ActorTable (first sheet) becomes the active / selected tab again
$powerLevelTable = $wb.Worksheets.Item("ActorPowerLevelTable")
$powerLevelTable.Range("F1").Select()

$actorTable.Activate()
